$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - "description metatag" -> "metatag", and its explanation text expanded
# to also mention the keywords metatag.
$ws.Range("B3").Value = "metatag"
$ws.Range("C3").Value = "The description tag is not helping in site ranking since 2009 according to Google but it is used to attract clicks when the website is searched on google, also  the keywords metatag is not usefull"

# Row 4 - explanation text tweaked ("a version" -> "a minified version", trailing period removed)
$ws.Range("C4").Value = "Linked css and js files are not minified, a minified version would reduce the size of the curent download"

# Category column capitalization / wording fixes
$ws.Range("A5").Value = "SEO/Accessibility"
$ws.Range("A6").Value = "SEO/Performance"
$ws.Range("A7").Value = "Performance/Accessibility"
$ws.Range("A8").Value = "Accessibility"
$ws.Range("A9").Value = "Accessibility"

# Update the sheet view: scroll so row 6 is at the top-left, and select B9.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("B9").Select()
